$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.35360442874736009,
    -0.0099999993021278044,
    -0.0089999992899087999,
    -0.011999999824686469,
    -0.0059999992980035444,
    -0.0059999992769519395,
    -0.019999999156699033,
    -0.019999999154191705,
    -0.0059999992721477824,
    -0.0059999992720349837,
    -0.0044999992847962744,
    -0.0059999992724217854,
    0.072810387280820343,
    -0.01199999922024908,
    -0.0059999992671206925,
    -0.025229882244500068,
    -0.0059999992598962493,
    -0.0089999992338700707,
    -0.040969755683968323,
    -0.0089999993036311565,
    -0.0089999993027145564,
    -0.0089999993021576685,
    -0.0089999992842164644,
    -0.041999998992732657,
    -0.041999998987352072,
    -0.0059999992750405795,
    -0.0059999992718879902,
    -0.005999999257495503,
    -0.011999999197486844,
    -0.019999999124447054,
    -0.014999999159709176,
    -0.021554311818968941,
    -0.0059999992348416242
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
